$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price report (2022-10-11, serial 44845) was added to the weekly
# series. It belongs at the top of this block (most recent date), so shift
# the existing rows 74:110 down by one to make room at row 74.
$ws.Rows("74:74").Insert()

# Populate the newly opened row 74 with the new record.
$ws.Range("A74").Value = 9
$ws.Range("B74").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C74").Value = "Metropolitana"
$ws.Range("D74").Value = 44845
$ws.Range("E74").Value = 13
$ws.Range("F74").Value = 100114002
$ws.Range("G74").Value = "Camote"
$ws.Range("H74").Value = "Sin especificar"
$ws.Range("I74").Value = "Primera"
$ws.Range("J74").Value = 600
$ws.Range("K74").Value = 18000
$ws.Range("L74").Value = 18000
$ws.Range("M74").Value = 18000
$ws.Range("N74").Value = "$/malla 18 kilos"
$ws.Range("O74").Value = "Perú"
$ws.Range("P74").Value = 1000
$ws.Range("Q74").Value = 18
$ws.Range("R74").Value = "Hortaliza"
